# huashu0_copy.xlsx -- "add flasksocketio evelent can run"
#
# 1. Re-home the frozen pane on Sheet1 so the visible viewport starts at
#    row 2 (just below the frozen header row) instead of row 219, while
#    leaving the current selection (C238) untouched.
# 2. On 工作表1, re-enter the volatile RAND() formula across C1:C18 as a
#    single range write so Excel collapses it into one shared formula
#    group (anchor C1 carries the formula text + ref span, C2:C18 just
#    reference the shared group) and gets fresh recalculated values.
# 3. Best-effort: nudge the workbook window's placement / minimized state
#    to match the refreshed window geometry recorded in the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet1: scroll the frozen pane back up to A2, keep selection C238 ---
$sheet1 = $wb.Worksheets.Item(1)
$sheet1.Activate()
$sheet1.Range("C238").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1

# --- 工作表1: rebuild C1:C18 as one shared RAND() formula ---
$sheet2 = $wb.Worksheets.Item(2)
$sheet2.Range("C1:C18").Formula = "=RAND()"

# --- Workbook window bookkeeping (best effort; harmless if unsupported) ---
try { $win.WindowState = -4140 } catch { }
try { $win.Left = 1880 } catch { }
try { $win.Top = 8280 } catch { }

$excel.CalculateFull()
